# Apply cryptocurrency price/volume updates to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.933.80"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.633.21"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "1.859.18"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "1.615.78"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₃0756"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "25.923.33"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.903"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "1.136.69"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.805"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").Value = "1.768.40"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  +2.90%  "
